$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data (B:G) for rows 2..11 after the update.
# Row 2 is brand-new data; rows 3..11 are the old rows 2..10 shifted down by one.
$data = @(
    @(0.01140061561852912, 0.1285420665309999, 0.03014188771060794, 0.1736141921347674, 0.1779866349240244, 19),
    @(0.03281260289867554, 0.712211698162994, 1.276354782257036, 1.12975872745336, 1.162021718847977, 18),
    @(-0.08894343397640261, 0.6850734014878143, 1.351404884049018, 1.162499412494053, 1.194764546449047, 17),
    @(0.01098408408208978, 0.5782217553777625, 0.6081029103588239, 0.7798095346678084, 0.805303924577184, 16),
    @(0.1235330352319591, 0.5509845820580629, 0.7140980823158125, 0.8450432428673769, 0.8653061007871674, 15),
    @(0.08484111210645341, 0.6180194413583581, 0.73762096509949, 0.858848627582003, 0.8869099907129777, 14),
    @(0.09014767153583389, 0.6854765477197997, 0.9025650063113468, 0.9500342132320009, 0.9843652458633849, 13),
    @(0.1184762789139503, 0.603487522271732, 0.7334425100585387, 0.856412581679262, 0.8858930048630905, 12),
    @(0.1655637249968011, 0.6554059886689477, 0.9061585844750577, 0.9519236232361595, 0.9831693451202929, 11),
    @(0.1561833212529775, 0.8228064401698413, 1.151759043685242, 1.073200374433983, 1.119208953704769, 10)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
